# Add two new worksheets ("Drinks" and "Games") with raw data + JSON-building
# formulas, matching the commit "Add json files for raw data".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Create the "Games" sheet first, then "Drinks" right after Sheet1 so the
# resulting tab order is Sheet1, Drinks, Games (matching workbook sheetId
# allocation: Games gets the lower internal id because it is created first).
# NOTE: worksheet object handles returned by .Add() are positional, not
# identity-bound -- inserting a second sheet at the same "After" position
# shifts the first one and its handle starts referring to the new sheet.
# So we re-fetch handles by name once both sheets exist.
$null = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$wb.Worksheets.Item(2).Name = "Games"
$null = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$wb.Worksheets.Item(2).Name = "Drinks"

$wsGames = $wb.Worksheets.Item("Games")
$wsDrinks = $wb.Worksheets.Item("Drinks")

# ---------------------------------------------------------------------
# Games sheet
# ---------------------------------------------------------------------
$wsGames.Range("E2").Value = "Followers"

$games = @(
    @(1, "Soccer", 70577),
    @(2, "Basketball", 20744),
    @(3, "Tennis", 10003),
    @(4, "Baseball", 23992),
    @(5, "Golf", 3456),
    @(6, "Running", 1205478),
    @(7, "Volleyball", 50944),
    @(8, "Badminton", 126087),
    @(9, "Swimming", 256843),
    @(10, "Boxing", 177435),
    @(11, "Table Tennis", 198436)
)

$row = 3
foreach ($g in $games) {
    $wsGames.Range("C$row").Value = $g[0]
    $wsGames.Range("D$row").Value = $g[1]
    $wsGames.Range("E$row").Value = $g[2]
    $row = $row + 1
}

# Center-aligned style for column C (id)
$wsGames.Range("C3:C13").HorizontalAlignment = -4108
$wsGames.Range("C3:C13").VerticalAlignment = -4108

# Center-aligned + bordered style for column D (game name) - reuse the
# existing border (thin, indexed color 64) already present on Sheet1 by
# copying its formats, then re-apply alignment.
$ws1.Range("B4").Copy()
$wsGames.Range("D3:D13").PasteSpecial(-4122)
$wsGames.Range("D3:D13").HorizontalAlignment = -4108
$wsGames.Range("D3:D13").VerticalAlignment = -4108

# Integer number format for column E (followers)
$wsGames.Range("E3:E13").NumberFormat = "0"

# JSON-builder formula in column F (shared across F3:F13)
$wsGames.Range("F3:F13").Formula = "=""{ 'id' = ""&C3&"", 'Game' : ""&D3&"", 'Followers' : ""&E3&""},"""

$wsGames.Columns.Item(4).ColumnWidth = 16.57
$wsGames.Columns.Item(5).ColumnWidth = 10.57

$wsGames.Range("E19").Select()

# ---------------------------------------------------------------------
# Drinks sheet
# ---------------------------------------------------------------------
$wsDrinks.Range("E3").Value = "Price"
$wsDrinks.Range("F3").Value = "Available Quantity"

$drinks = @(
    @(1, "Wine", 250, 24456),
    @(2, "Coffee", 50, 3456),
    @(3, "Lemonade", 30, 2456),
    @(4, "Iced Tea", 40, 8743),
    @(5, "Hot Chocolate", 60, 7963),
    @(6, "Juice", 30, 6437),
    @(7, "Milkshake", 40, 8453),
    @(8, "Water", 20, 7357),
    @(9, "Tea", 10, 8442),
    @(10, "Beer", 140, 7457)
)

$row = 4
foreach ($d in $drinks) {
    $wsDrinks.Range("C$row").Value = $d[0]
    $wsDrinks.Range("D$row").Value = $d[1]
    $wsDrinks.Range("E$row").Value = $d[2]
    $wsDrinks.Range("F$row").Value = $d[3]
    $row = $row + 1
}

# Bordered style for column D (drink name) - reuse existing Sheet1 border
# style (s=2 in styles.xml) by copying formats from a cell that already has it.
$ws1.Range("B4").Copy()
$wsDrinks.Range("D4:D13").PasteSpecial(-4122)

# Integer number format for column F (availability)
$wsDrinks.Range("F4:F13").NumberFormat = "0"

# JSON-builder formula in column I (shared across I4:I13)
$wsDrinks.Range("I4:I13").Formula = "="" { 'id' : ""&C4&"", 'Drink' : '""&D4&""', 'Availability' : ""&F4&"", 'Price' : ""&E4&"" },"""

$wsDrinks.Range("O8").Select()

# ---------------------------------------------------------------------
# Sheet1 view tweak
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("J18").Select()
